$wb = $excel.ActiveWorkbook

# --- Sheet: "32bit Fixed Point (arctan)" (4th sheet) ---
$ws = $wb.Worksheets.Item(4)

# Step 1: H11 value fixed first (matches shared-string insertion order observed in the diff)
$ws.Range("H11").Value = "56/62"

# Step 2: I2:I13 bulk updated to "Every Sample"
$ws.Range("I2:I13").Value = "Every Sample"

# Step 3: G11, G12, G13 updated
$ws.Range("G11").Value = "56/58"
$ws.Range("G12").Value = "6/-2"
$ws.Range("G13").Value = "6/-6"

# Step 4: H12, H13 updated (H12 additionally gets a date-like number format applied)
$ws.Range("H12").Value = "1/1--"
$ws.Range("H13").Value = "2/-4"

# View state change on this sheet
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("A12:J13").Select()

Write-Host "done"
